$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of data for subject S5 / Marijn, session 2013-06-06-marijn
$rows = @(
    @("S5","Marijn","2013-06-06-marijn","2013-06-06-15-10-23","hybrid-15Hz",15,1),
    @("S5","Marijn","2013-06-06-marijn","2013-06-06-15-18-33","ssvep-10Hz",10,0),
    @("S5","Marijn","2013-06-06-marijn","2013-06-06-15-24-10","hybrid-12Hz",12,1),
    @("S5","Marijn","2013-06-06-marijn","2013-06-06-15-29-56","hybrid-10Hz",10,1),
    @("S5","Marijn","2013-06-06-marijn","2013-06-06-15-47-50","ssvep-15Hz",15,0),
    @("S5","Marijn","2013-06-06-marijn","2013-06-06-15-53-20","ssvep-12Hz",12,0),
    @("S5","Marijn","2013-06-06-marijn","2013-06-06-15-58-51","hybrid-8-57Hz",8.57,1),
    @("S5","Marijn","2013-06-06-marijn","2013-06-06-16-04-46","ssvep-8-57Hz",8.57,0)
)

$startRow = 34
$lastExistingRow = 33

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]

    # Date: reuse the numeric date serial (2013-06-06) and copy the
    # existing date cell's style so formatting matches the rest of the column.
    $ws.Cells.Item($r, 3).Value = 41431
    $ws.Cells.Item($lastExistingRow, 3).Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)

    $ws.Cells.Item($r, 4).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[3]
    $ws.Cells.Item($r, 6).Value = $data[4]
    $ws.Cells.Item($r, 7).Value = $data[5]
    $ws.Cells.Item($r, 8).Value = $data[6]
}

$excel.CutCopyMode = 0

$ws.Range("A35:D41").Select()
